$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert two new columns before column D (two new quarterly reporting periods),
# shifting the existing historical data from D:K to F:M
$ws.Range("D:E").EntireColumn.Insert()

# Copy number formatting from column F (the original column D, now shifted right)
# into the new D:E columns so each row keeps its correct style (date rows, numeric rows, etc.)
$ws.Range("F7:F102").Copy()
$ws.Range("D7:E102").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Populate the two new quarter columns with the latest reported figures
$ws.Range("D7").Value = 43490
$ws.Range("E7").Value = 43399
$ws.Range("D8").Value = 7546000
$ws.Range("E8").Value = 7481000
$ws.Range("D9").Value = 2244000
$ws.Range("E9").Value = 2181000
$ws.Range("D10").Value = 5302000
$ws.Range("E10").Value = 5300000
$ws.Range("D12").Value = 561000
$ws.Range("E12").Value = 590000
$ws.Range("D13:E13").Value = 0
$ws.Range("D14").Value = 146000
$ws.Range("E14").Value = 81000
$ws.Range("D15").Value = 436000
$ws.Range("E15").Value = 445000
$ws.Range("D17").Value = 6004000
$ws.Range("E17").Value = 5937000
$ws.Range("D18").Value = 1542000
$ws.Range("E18").Value = 1544000
$ws.Range("D20").Value = 71000
$ws.Range("E20").Value = 52000
$ws.Range("D21").Value = 2288000
$ws.Range("E21").Value = 2247000
$ws.Range("D22").Value = 243000
$ws.Range("E22").Value = 241000
$ws.Range("D23").Value = 1370000
$ws.Range("E23").Value = 1355000
$ws.Range("D24").Value = 111000
$ws.Range("E24").Value = 198000
$ws.Range("D25:E25").Value = 0
$ws.Range("D26").Value = 1259000
$ws.Range("E26").Value = 1157000
$ws.Range("D27").Value = 1257000
$ws.Range("E27").Value = 1152000
$ws.Range("D28:E28").Value = 0
$ws.Range("D29").Value = 12000
$ws.Range("E29").Value = -37000
$ws.Range("D30:E30").Value = 0
$ws.Range("D31:E31").Value = 0
$ws.Range("D32").Value = -71000
$ws.Range("E32").Value = -52000
$ws.Range("D33").Value = 1269000
$ws.Range("E33").Value = 1115000
$ws.Range("D34:E34").Value = 0
$ws.Range("D35").Value = 1269000
$ws.Range("E35").Value = 1115000
$ws.Range("D38").Value = 43490
$ws.Range("E38").Value = 43399
$ws.Range("D41").Value = 3703000
$ws.Range("E41").Value = 3911000
$ws.Range("D42").Value = 5439000
$ws.Range("E42").Value = 6222000
$ws.Range("D43").Value = 5854000
$ws.Range("E43").Value = 5743000
$ws.Range("D44").Value = 3866000
$ws.Range("E44").Value = 3763000
$ws.Range("D45").Value = 2015000
$ws.Range("E45").Value = 2014000
$ws.Range("D46").Value = 20877000
$ws.Range("E46").Value = 21653000
$ws.Range("D47").Value = 330000
$ws.Range("E47").Value = 511000
$ws.Range("D48").Value = 4593000
$ws.Range("E48").Value = 4536000
$ws.Range("D49").Value = 60838000
$ws.Range("E49").Value = 59424000
$ws.Range("D50:E50").Value = 0
$ws.Range("D51:E51").Value = 0
$ws.Range("D52").Value = 2092000
$ws.Range("E52").Value = 2026000
$ws.Range("D53:E53").Value = 0
$ws.Range("D54").Value = 88730000
$ws.Range("E54").Value = 88150000
$ws.Range("D57").Value = 1706000
$ws.Range("E57").Value = 1742000
$ws.Range("D58").Value = 1356000
$ws.Range("E58").Value = 1343000
$ws.Range("D59").Value = 5791000
$ws.Range("E59").Value = 5378000
$ws.Range("D60").Value = 8853000
$ws.Range("E60").Value = 8463000
$ws.Range("D61").Value = 23674000
$ws.Range("E61").Value = 23673000
$ws.Range("D62").Value = 6262000
$ws.Range("E62").Value = 6300000
$ws.Range("D63:E63").Value = 0
$ws.Range("D64:E64").Value = 0
$ws.Range("D65:E65").Value = 0
$ws.Range("D66").Value = 38901000
$ws.Range("E66").Value = 38543000
$ws.Range("D68:E68").Value = 0
$ws.Range("D69:E69").Value = 0
$ws.Range("D70:E70").Value = 0
$ws.Range("D71:E71").Value = 0
$ws.Range("D72").Value = 25769000
$ws.Range("E72").Value = 25171000
$ws.Range("D73:E73").Value = 0
$ws.Range("D74:E74").Value = 0
$ws.Range("D75:E75").Value = 0
$ws.Range("D76").Value = 49829000
$ws.Range("E76").Value = 49607000
$ws.Range("D77:E77").Value = 0
$ws.Range("D80").Value = 43490
$ws.Range("E80").Value = 43399
$ws.Range("D81").Value = 1269000
$ws.Range("E81").Value = 1115000
$ws.Range("D83").Value = 675000
$ws.Range("E83").Value = 651000
$ws.Range("D84:E84").Value = 0
$ws.Range("D85:E85").Value = 0
$ws.Range("D86:E86").Value = 0
$ws.Range("D87:E87").Value = 0
$ws.Range("D88:E88").Value = 0
$ws.Range("D89").Value = 2055000
$ws.Range("E89").Value = 1163000
$ws.Range("D91").Value = -302000
$ws.Range("E91").Value = -206000
$ws.Range("D92:E92").Value = 0
$ws.Range("D93:E93").Value = 0
$ws.Range("D94").Value = -1009000
$ws.Range("E94").Value = 121000
$ws.Range("D96").Value = -671000
$ws.Range("E96").Value = -674000
$ws.Range("D97:E97").Value = 0
$ws.Range("D98:E98").Value = 0
$ws.Range("D99:E99").Value = 0
$ws.Range("D100").Value = -1268000
$ws.Range("E100").Value = -1730000
$ws.Range("D101").Value = 14000
$ws.Range("E101").Value = -23000
$ws.Range("D102").Value = -208000
$ws.Range("E102").Value = -469000

# Apply minor restatements to previously reported historical quarters
$ws.Range("H9").Value = 2193000
$ws.Range("I9").Value = 2116000
$ws.Range("J9").Value = 2347000
$ws.Range("H10").Value = 5176000
$ws.Range("I10").Value = 4934000
$ws.Range("J10").Value = 5043000
$ws.Range("H12").Value = 559000
$ws.Range("I12").Value = 556000
$ws.Range("F14").Value = 252000
$ws.Range("H14").Value = 105000
$ws.Range("I14").Value = -514000
$ws.Range("J14").Value = 106000
$ws.Range("H17").Value = 5933000
$ws.Range("I17").Value = 5156000
$ws.Range("H18").Value = 1436000
$ws.Range("I18").Value = 1894000
$ws.Range("H20").Value = -139000
$ws.Range("I20").Value = 107000
$ws.Range("F24").Value = 153000
$ws.Range("F26").Value = 1027000
$ws.Range("F27").Value = 1025000
$ws.Range("F29").Value = 50000
$ws.Range("H32").Value = 139000
$ws.Range("I32").Value = -107000
